# The workbook's single sheet stores forest-ministry login rows. This
# update renames the header columns to match the new DB schema
# (ministryName/username -> ministryID/locationID) and removes the
# leftover "Forest Ministry" label that had been duplicated into A2,
# leaving that cell blank while every other row (and the whole B column)
# keeps its original fm0xx values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header names
$ws.Range("A1").Value = "ministryID"
$ws.Range("B1").Value = "locationID"

# A2 no longer carries the "Forest Ministry" text
$ws.Range("A2").ClearContents()

# Column B's best-fit width grows slightly now that "locationID" (header)
# is wider than the old "username" header.
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666

# Active selection ends up on A2 after the edit
$ws.Range("A2").Select()
